$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text field fixes (comma -> period typos in Razon social / Nombre Fantasia) ---
foreach ($addr in @("E34","F34","E60","F60","E69","F69")) {
  $ws.Range($addr).Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
}
foreach ($addr in @("E56","F56")) {
  $ws.Range($addr).Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
}
foreach ($addr in @("E59")) {
  $ws.Range($addr).Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
}
foreach ($addr in @("E61","E123")) {
  $ws.Range($addr).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
}

# --- Importe column (H2:H160): Spanish-formatted numeric text -> plain dotted-decimal text ---
$importeRange = $ws.Range("H2:H160")
$importeRange.NumberFormat = "@"

$importeValues = @{
  2 = "654.00"
  3 = "12907.07"
  4 = "1600.00"
  5 = "1000.00"
  6 = "750.00"
  7 = "991.00"
  8 = "105.00"
  9 = "3324.50"
  10 = "223.00"
  11 = "7299.94"
  12 = "15558.53"
  13 = "18336.90"
  14 = "49851.05"
  15 = "6785.38"
  16 = "2672.47"
  17 = "19248.45"
  18 = "5932.08"
  19 = "86.00"
  20 = "1160.90"
  21 = "2.00"
  22 = "8013.12"
  23 = "650.00"
  24 = "129.28"
  25 = "38.00"
  26 = "749.00"
  27 = "109.00"
  28 = "181.88"
  29 = "2241.89"
  30 = "6721.02"
  31 = "69.06"
  32 = "2381.57"
  33 = "7.32"
  34 = "55.45"
  35 = "1845.08"
  36 = "30.00"
  37 = "1416.04"
  38 = "67.00"
  39 = "256.00"
  40 = "5907.40"
  41 = "16.94"
  42 = "58.00"
  43 = "6202.09"
  44 = "85.00"
  45 = "3427.73"
  46 = "3669.02"
  47 = "925.00"
  48 = "5.00"
  49 = "381.97"
  50 = "315.00"
  51 = "402.53"
  52 = "417.45"
  53 = "800.00"
  54 = "4190.00"
  55 = "4573.80"
  56 = "70.00"
  57 = "11094.00"
  58 = "1357.50"
  59 = "625.00"
  60 = "3077.19"
  61 = "1907.00"
  62 = "11.00"
  63 = "71.00"
  64 = "212532.00"
  65 = "156.24"
  66 = "2118.15"
  67 = "1.47"
  68 = "2408.93"
  69 = "22.41"
  70 = "107.74"
  71 = "7.98"
  72 = "205.84"
  73 = "1023.90"
  74 = "509.88"
  75 = "8.70"
  76 = "39.30"
  77 = "600.00"
  78 = "4706.30"
  79 = "1457.55"
  80 = "72.00"
  81 = "165.75"
  82 = "50.63"
  83 = "725.00"
  84 = "175.50"
  85 = "4399.46"
  86 = "439.00"
  87 = "52188.00"
  88 = "12520.00"
  89 = "630.00"
  90 = "2500.00"
  91 = "300.00"
  92 = "110.00"
  93 = "250.00"
  94 = "105.12"
  95 = "132.30"
  96 = "64.00"
  97 = "693.50"
  98 = "322.50"
  99 = "28.34"
  100 = "207.30"
  101 = "17530.00"
  102 = "41098.86"
  103 = "1594.00"
  104 = "400.00"
  105 = "1000.00"
  106 = "300.00"
  107 = "100.00"
  108 = "1500.00"
  109 = "5755.27"
  110 = "170.00"
  111 = "450.00"
  112 = "600.00"
  113 = "200.00"
  114 = "4137.14"
  115 = "660.00"
  116 = "9111.30"
  117 = "500.00"
  118 = "800.00"
  119 = "240.00"
  120 = "3000.00"
  121 = "744.01"
  122 = "206.00"
  123 = "201.00"
  124 = "186.00"
  125 = "20.00"
  126 = "1903.50"
  127 = "290.00"
  128 = "200.00"
  129 = "230.00"
  130 = "50.00"
  131 = "190.21"
  132 = "8398.00"
  133 = "699.17"
  134 = "634.00"
  135 = "628.00"
  136 = "20860.00"
  137 = "10.40"
  138 = "149.07"
  139 = "3225.43"
  140 = "1698.67"
  141 = "2920.00"
  142 = "3473.37"
  143 = "3900.00"
  144 = "283.40"
  145 = "5.27"
  146 = "112.30"
  147 = "1602.00"
  148 = "9427.00"
  149 = "407.73"
  150 = "5790.46"
  151 = "7946.31"
  152 = "1800.00"
  153 = "572003.68"
  154 = "11495.00"
  155 = "4500.00"
  156 = "875.60"
  157 = "4089.40"
  158 = "750.00"
  159 = "4100.00"
  160 = "276.00"
}

foreach ($r in $importeValues.Keys) {
  $ws.Cells.Item($r, 8).Value = $importeValues[$r]
}

$importeRange.Style = "Normal"
